$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (20) of forecast data, mirroring the existing rows' layout/formatting
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value = 45986

$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = -0.08656168856399082
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = -0.02867614772544824
